$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.063706102998296
$ws.Cells.Item(2, 4).Value = 1.071065694767837
$ws.Cells.Item(2, 5).Value = 1.077239540357198
$ws.Cells.Item(2, 6).Value = 1.084178910349178
$ws.Cells.Item(2, 9).Value = 1.059266692982902
$ws.Cells.Item(2, 10).Value = 1.068670102219458
$ws.Cells.Item(2, 11).Value = 1.073763859730122
$ws.Cells.Item(2, 12).Value = 1.079921362481024
$ws.Cells.Item(2, 13).Value = 1.086842603085771
$ws.Cells.Item(2, 14).Value = 1.070187736342784
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.064695555733248
$ws.Cells.Item(3, 4).Value = 1.071877746886015
$ws.Cells.Item(3, 5).Value = 1.078182794255318
$ws.Cells.Item(3, 6).Value = 1.085123451853068
$ws.Cells.Item(3, 9).Value = 1.059595675439569
$ws.Cells.Item(3, 10).Value = 1.069314293348623
$ws.Cells.Item(3, 11).Value = 1.074392404809425
$ws.Cells.Item(3, 12).Value = 1.080681965098084
$ws.Cells.Item(3, 13).Value = 1.087605798028291
$ws.Cells.Item(3, 14).Value = 1.070832842297242
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.065336284078112
$ws.Cells.Item(4, 4).Value = 1.072403592756781
$ws.Cells.Item(4, 5).Value = 1.078793925081053
$ws.Cells.Item(4, 6).Value = 1.08573538122412
$ws.Cells.Item(4, 9).Value = 1.059807540231874
$ws.Cells.Item(4, 10).Value = 1.069730964947364
$ws.Cells.Item(4, 11).Value = 1.074798853240231
$ws.Cells.Item(4, 12).Value = 1.081174277035511
$ws.Cells.Item(4, 13).Value = 1.088099751964263
$ws.Cells.Item(4, 14).Value = 1.071250105617446
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.065605761885952
$ws.Cells.Item(5, 4).Value = 1.072624751489383
$ws.Cells.Item(5, 5).Value = 1.079051030647368
$ws.Cells.Item(5, 6).Value = 1.085992814176704
$ws.Cells.Item(5, 9).Value = 1.059896366156977
$ws.Cells.Item(5, 10).Value = 1.069906093814685
$ws.Cells.Item(5, 11).Value = 1.074969660510102
$ws.Cells.Item(5, 12).Value = 1.081381280205099
$ws.Cells.Item(5, 13).Value = 1.088307436921399
$ws.Cells.Item(5, 14).Value = 1.071425483187846
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.065651015152031
$ws.Cells.Item(6, 4).Value = 1.072661890441155
$ws.Cells.Item(6, 5).Value = 1.079094210686544
$ws.Cells.Item(6, 6).Value = 1.086036048696476
$ws.Cells.Item(6, 9).Value = 1.059911266212887
$ws.Cells.Item(6, 10).Value = 1.069935496360258
$ws.Cells.Item(6, 11).Value = 1.074998336044761
$ws.Cells.Item(6, 12).Value = 1.081416038986179
$ws.Cells.Item(6, 13).Value = 1.088342309677738
$ws.Cells.Item(6, 14).Value = 1.071454927488407
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.065339884400481
$ws.Cells.Item(7, 4).Value = 1.072406547526107
$ws.Cells.Item(7, 5).Value = 1.07879735980815
$ws.Cells.Item(7, 6).Value = 1.085738820358511
$ws.Cells.Item(7, 9).Value = 1.059808728079891
$ws.Cells.Item(7, 10).Value = 1.06973330518555
$ws.Cells.Item(7, 11).Value = 1.074801135826509
$ws.Cells.Item(7, 12).Value = 1.081177042884759
$ws.Cells.Item(7, 13).Value = 1.088102526957063
$ws.Cells.Item(7, 14).Value = 1.071252449179038
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.064040392015938
$ws.Cells.Item(8, 4).Value = 1.071340049336207
$ws.Cells.Item(8, 5).Value = 1.077558154712495
$ws.Cells.Item(8, 6).Value = 1.084497967005323
$ws.Cells.Item(8, 9).Value = 1.059378082726288
$ws.Cells.Item(8, 10).Value = 1.068887842914274
$ws.Cells.Item(8, 11).Value = 1.073976333266351
$ws.Cells.Item(8, 12).Value = 1.080178379813397
$ws.Cells.Item(8, 13).Value = 1.087100503819464
$ws.Cells.Item(8, 14).Value = 1.070405786254361
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.0617542799393
$ws.Cells.Item(9, 4).Value = 1.069463810870434
$ws.Cells.Item(9, 5).Value = 1.075380561082846
$ws.Cells.Item(9, 6).Value = 1.082317205259714
$ws.Cells.Item(9, 9).Value = 1.058611526576389
$ws.Cells.Item(9, 10).Value = 1.067396817009
$ws.Cells.Item(9, 11).Value = 1.072520955257938
$ws.Cells.Item(9, 12).Value = 1.078419808273037
$ws.Cells.Item(9, 13).Value = 1.085335741612817
$ws.Cells.Item(9, 14).Value = 1.068912642921279
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.060232779039163
$ws.Cells.Item(10, 4).Value = 1.068215115811802
$ws.Cells.Item(10, 5).Value = 1.073932963664907
$ws.Cells.Item(10, 6).Value = 1.080867320497697
$ws.Cells.Item(10, 9).Value = 1.058095339273835
$ws.Cells.Item(10, 10).Value = 1.066402031663775
$ws.Cells.Item(10, 11).Value = 1.071549428906875
$ws.Cells.Item(10, 12).Value = 1.077248290625739
$ws.Cells.Item(10, 13).Value = 1.084159917342066
$ws.Cells.Item(10, 14).Value = 1.06791644486676
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.059574572042825
$ws.Cells.Item(11, 4).Value = 1.067674936475974
$ws.Cells.Item(11, 5).Value = 1.073307131962776
$ws.Cells.Item(11, 6).Value = 1.080240457207683
$ws.Cells.Item(11, 9).Value = 1.057870609074047
$ws.Cells.Item(11, 10).Value = 1.065971106204396
$ws.Cells.Item(11, 11).Value = 1.071128455371004
$ws.Cells.Item(11, 12).Value = 1.076741226879817
$ws.Cells.Item(11, 13).Value = 1.083650946602639
$ws.Cells.Item(11, 14).Value = 1.067484907443801
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.059330177421941
$ws.Cells.Item(12, 4).Value = 1.067474368118723
$ws.Cells.Item(12, 5).Value = 1.073074819656336
$ws.Cells.Item(12, 6).Value = 1.08000775557549
$ws.Cells.Item(12, 9).Value = 1.057786951775648
$ws.Cells.Item(12, 10).Value = 1.065811015475585
$ws.Cells.Item(12, 11).Value = 1.07097204334436
$ws.Cells.Item(12, 12).Value = 1.076552913434086
$ws.Cells.Item(12, 13).Value = 1.083461918507139
$ws.Cells.Item(12, 14).Value = 1.067324589367793
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.059382596689014
$ws.Cells.Item(13, 4).Value = 1.067517387151346
$ws.Cells.Item(13, 5).Value = 1.073124644654126
$ws.Cells.Item(13, 6).Value = 1.080057664362635
$ws.Cells.Item(13, 9).Value = 1.057804904801234
$ws.Cells.Item(13, 10).Value = 1.065845356637102
$ws.Cells.Item(13, 11).Value = 1.071005596220334
$ws.Cells.Item(13, 12).Value = 1.076593305817651
$ws.Cells.Item(13, 13).Value = 1.083502464470417
$ws.Cells.Item(13, 14).Value = 1.067358979297698
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.059554368410916
$ws.Cells.Item(14, 4).Value = 1.067658355829356
$ws.Cells.Item(14, 5).Value = 1.073287925899961
$ws.Cells.Item(14, 6).Value = 1.080221219088299
$ws.Cells.Item(14, 9).Value = 1.057863697655922
$ws.Cells.Item(14, 10).Value = 1.065957873577355
$ws.Cells.Item(14, 11).Value = 1.071115527191203
$ws.Cells.Item(14, 12).Value = 1.076725660161341
$ws.Cells.Item(14, 13).Value = 1.083635320939475
$ws.Cells.Item(14, 14).Value = 1.067471656024912
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.059660214950512
$ws.Cells.Item(15, 4).Value = 1.067745221692734
$ws.Cells.Item(15, 5).Value = 1.073388548701014
$ws.Cells.Item(15, 6).Value = 1.080322009575655
$ws.Cells.Item(15, 9).Value = 1.057899897705076
$ws.Cells.Item(15, 10).Value = 1.066027195577798
$ws.Cells.Item(15, 11).Value = 1.07118325352205
$ws.Cells.Item(15, 12).Value = 1.076807212388006
$ws.Cells.Item(15, 13).Value = 1.083717181707312
$ws.Cells.Item(15, 14).Value = 1.067541076470547
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.060276475041286
$ws.Cells.Item(16, 4).Value = 1.068250976655442
$ws.Cells.Item(16, 5).Value = 1.073974518952774
$ws.Cells.Item(16, 6).Value = 1.080908943384574
$ws.Cells.Item(16, 9).Value = 1.058110228259889
$ws.Cells.Item(16, 10).Value = 1.066430627104741
$ws.Cells.Item(16, 11).Value = 1.071577361360588
$ws.Cells.Item(16, 12).Value = 1.077281947300431
$ws.Cells.Item(16, 13).Value = 1.084193699691338
$ws.Cells.Item(16, 14).Value = 1.067945080916532
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.060663203291303
$ws.Cells.Item(17, 4).Value = 1.068568361731965
$ws.Cells.Item(17, 5).Value = 1.074342347767835
$ws.Cells.Item(17, 6).Value = 1.081277365620322
$ws.Cells.Item(17, 9).Value = 1.058241837364877
$ws.Cells.Item(17, 10).Value = 1.066683642220804
$ws.Cells.Item(17, 11).Value = 1.071824496048227
$ws.Cells.Item(17, 12).Value = 1.077579793118569
$ws.Cells.Item(17, 13).Value = 1.084492652708552
$ws.Cells.Item(17, 14).Value = 1.068198455343082
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.060888834430054
$ws.Cells.Item(18, 4).Value = 1.068753536584203
$ws.Cells.Item(18, 5).Value = 1.074556991302002
$ws.Cells.Item(18, 6).Value = 1.081492351313769
$ws.Cells.Item(18, 9).Value = 1.058318485209341
$ws.Cells.Item(18, 10).Value = 1.066831204544874
$ws.Cells.Item(18, 11).Value = 1.071968616896691
$ws.Cells.Item(18, 12).Value = 1.077753541856097
$ws.Cells.Item(18, 13).Value = 1.084667043183005
$ws.Cells.Item(18, 14).Value = 1.068346227222577
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.060965778820275
$ws.Cells.Item(19, 4).Value = 1.068816684787327
$ws.Cells.Item(19, 5).Value = 1.074630195316411
$ws.Cells.Item(19, 6).Value = 1.08156567132071
$ws.Cells.Item(19, 9).Value = 1.058344600196558
$ws.Cells.Item(19, 10).Value = 1.06688151655893
$ws.Cells.Item(19, 11).Value = 1.072017753500802
$ws.Cells.Item(19, 12).Value = 1.077812789103667
$ws.Cells.Item(19, 13).Value = 1.084726508554356
$ws.Cells.Item(19, 14).Value = 1.068396610685464
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.060621704892779
$ws.Cells.Item(20, 4).Value = 1.06853430418828
$ws.Cells.Item(20, 5).Value = 1.074302873361223
$ws.Cells.Item(20, 6).Value = 1.081237827955979
$ws.Cells.Item(20, 9).Value = 1.058227729105819
$ws.Cells.Item(20, 10).Value = 1.066656497857994
$ws.Cells.Item(20, 11).Value = 1.071797983791072
$ws.Cells.Item(20, 12).Value = 1.07754783497377
$ws.Cells.Item(20, 13).Value = 1.084460576192887
$ws.Cells.Item(20, 14).Value = 1.068171272432163
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.059503783353893
$ws.Cells.Item(21, 4).Value = 1.067616841896964
$ws.Cells.Item(21, 5).Value = 1.073239839508365
$ws.Cells.Item(21, 6).Value = 1.080173052333149
$ws.Cells.Item(21, 9).Value = 1.057846389659758
$ws.Cells.Item(21, 10).Value = 1.065924740849818
$ws.Cells.Item(21, 11).Value = 1.071083156460524
$ws.Cells.Item(21, 12).Value = 1.076686684198731
$ws.Cells.Item(21, 13).Value = 1.083596197287318
$ws.Cells.Item(21, 14).Value = 1.067438476245101
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.058801438757596
$ws.Cells.Item(22, 4).Value = 1.067040450501911
$ws.Cells.Item(22, 5).Value = 1.072572333162019
$ws.Cells.Item(22, 6).Value = 1.079504415313049
$ws.Cells.Item(22, 9).Value = 1.05760557077317
$ws.Cells.Item(22, 10).Value = 1.065464507122411
$ws.Cells.Item(22, 11).Value = 1.070633463638554
$ws.Cells.Item(22, 12).Value = 1.076145433620197
$ws.Cells.Item(22, 13).Value = 1.083052880577762
$ws.Cells.Item(22, 14).Value = 1.066977588933012
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.059173713783215
$ws.Cells.Item(23, 4).Value = 1.067345963126939
$ws.Cells.Item(23, 5).Value = 1.072926108528389
$ws.Cells.Item(23, 6).Value = 1.079858793429524
$ws.Cells.Item(23, 9).Value = 1.057733333333393
$ws.Cells.Item(23, 10).Value = 1.065708499609137
$ws.Cells.Item(23, 11).Value = 1.070871877981796
$ws.Cells.Item(23, 12).Value = 1.076432342571836
$ws.Cells.Item(23, 13).Value = 1.083340888271451
$ws.Cells.Item(23, 14).Value = 1.067221927917055
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.060640456050211
$ws.Cells.Item(24, 4).Value = 1.068549693173404
$ws.Cells.Item(24, 5).Value = 1.074320709852316
$ws.Cells.Item(24, 6).Value = 1.081255693042798
$ws.Cells.Item(24, 9).Value = 1.058234104383939
$ws.Cells.Item(24, 10).Value = 1.066668763280049
$ws.Cells.Item(24, 11).Value = 1.071809963627397
$ws.Cells.Item(24, 12).Value = 1.077562275421082
$ws.Cells.Item(24, 13).Value = 1.08447507013949
$ws.Cells.Item(24, 14).Value = 1.068183555272524
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.062344844453266
$ws.Cells.Item(25, 4).Value = 1.069948492947683
$ws.Cells.Item(25, 5).Value = 1.075942797983651
$ws.Cells.Item(25, 6).Value = 1.082880292233347
$ws.Cells.Item(25, 9).Value = 1.058810609614231
$ws.Cells.Item(25, 10).Value = 1.067782421827729
$ws.Cells.Item(25, 11).Value = 1.072897433678209
$ws.Cells.Item(25, 12).Value = 1.078874293124783
$ws.Cells.Item(25, 13).Value = 1.08579185867825
$ws.Cells.Item(25, 14).Value = 1.06929879534308

Write-Host "Updated vm_pu values for case with 380 kV"
